# Add a space before ":" in the statut_name values for rows 2-4 (column B)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "1 : résultats postés ou publiés dans les 12 mois"
$ws.Range("B3").Value = "4 : pas de résultats postés ni publiés"
$ws.Range("B4").Value = "4 : pas de résultats postés ni publiés"
